# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet and
# moves the special "latest day" date formatting down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (21) currently carries the "latest row" date
# formatting (YYYY-MM-DD, no time). Since a new row is being appended,
# row 21 reverts to the regular date formatting used by all the other
# historical rows (YYYY-MM-DD HH:MM:SS, same as row 20).
$ws.Range("A21").NumberFormat = $ws.Range("A20").NumberFormat

# Append the new day's data as row 22.
$ws.Range("A22").Value = 45762
$ws.Range("B22").Value = 87
$ws.Range("C22").Value = 87
$ws.Range("D22").Value = 85

# Row 22 is now the newest/last row, so it gets the special date-only
# format that row 21 used to have.
$ws.Range("A22").NumberFormat = "YYYY-MM-DD"
